$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10003209
$ws.Range("I32").Value = 2163
$ws.Range("K32").Value = 2163
$ws.Range("M32").Value = -1837

$ws.Range("H39").Value = 2325.4285
$ws.Range("J39").Value = 3825
$ws.Range("L39").Value = 11475
$ws.Range("N39").Value = -12067

$ws.Range("H70").Value = 2244
$ws.Range("J70").Value = 2449.5
$ws.Range("L70").Value = 7348.5
$ws.Range("N70").Value = -7888.5

$ws.Range("H73").Value = 2244
$ws.Range("J73").Value = 2449.5
$ws.Range("L73").Value = 7348.5
$ws.Range("N73").Value = -9220.5

$ws.Range("H76").Value = 7596.75
$ws.Range("I76").Value = 7388
$ws.Range("K76").Value = 7388
$ws.Range("M76").Value = -7073

$ws.Range("H79").Value = 7596.75
$ws.Range("I79").Value = 7388
$ws.Range("K79").Value = 7388
$ws.Range("M79").Value = -6296

$ws.Range("H100").Value = 7346.294
$ws.Range("I100").Value = 6085.25
$ws.Range("K100").Value = 6085.25
$ws.Range("M100").Value = -5544.25

$ws.Range("H103").Value = 1281.8334
$ws.Range("I103").Value = 1098.3334
$ws.Range("J103").Value = 1465.3334
$ws.Range("K103").Value = 3295.0002
$ws.Range("L103").Value = 4396.0002
$ws.Range("M103").Value = -2709.0002
$ws.Range("N103").Value = -5568.0002

$ws.Range("H106").Value = 8167.5
$ws.Range("I106").Value = 8000.909
$ws.Range("K106").Value = 8000.909
$ws.Range("M106").Value = -7369.909

$ws.Range("H113").Value = 8671.888999999999
$ws.Range("I113").Value = 6775
$ws.Range("K113").Value = 6775
$ws.Range("M113").Value = -3521

$ws.Range("H130").Value = 32500
$ws.Range("J130").Value = 32500
$ws.Range("L130").Value = 32500
$ws.Range("N130").Value = -42540

$ws.Range("H132").Value = 8739.52
$ws.Range("I132").Value = 1545.9841
$ws.Range("K132").Value = 4637.9523
$ws.Range("M132").Value = -2107.9523

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2264.923
$ws.Range("I61").Value = 2057.8333
$ws.Range("K61").Value = 2057.8333
$ws.Range("M61").Value = -1845.8333

$ws.Range("H122").Value = 3783.3
$ws.Range("I122").Value = 2978.4211
$ws.Range("K122").Value = 8935.263300000001
$ws.Range("M122").Value = -6485.263300000001

$ws.Range("H136").Value = 2264.923
$ws.Range("I136").Value = 2057.8333
$ws.Range("K136").Value = 6173.499899999999
$ws.Range("M136").Value = -3623.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4583.457
$ws.Range("I107").Value = 3071.7856
$ws.Range("J107").Value = 10630.143
$ws.Range("K107").Value = 3071.7856
$ws.Range("L107").Value = 10630.143
$ws.Range("M107").Value = -1151.7856
$ws.Range("N107").Value = -14470.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 36666.668
$ws.Range("J9").Value = 36666.668
$ws.Range("L9").Value = 36666.668
$ws.Range("N9").Value = -37002.668

$ws.Range("H31").Value = 1539.7428
$ws.Range("I31").Value = 1400.5667
$ws.Range("J31").Value = 2374.8
$ws.Range("K31").Value = 1400.5667
$ws.Range("L31").Value = 2374.8
$ws.Range("M31").Value = -1105.5667
$ws.Range("N31").Value = -2964.8

$ws.Range("H34").Value = 1539.7428
$ws.Range("I34").Value = 1400.5667
$ws.Range("J34").Value = 2374.8
$ws.Range("K34").Value = 1400.5667
$ws.Range("L34").Value = 2374.8
$ws.Range("M34").Value = -1198.5667
$ws.Range("N34").Value = -2778.8

$ws.Range("H58").Value = 1910.4445
$ws.Range("I58").Value = 1242
$ws.Range("K58").Value = 1242
$ws.Range("M58").Value = -1039

$ws.Range("H107").Value = 7591.533
$ws.Range("I107").Value = 859.64703
$ws.Range("K107").Value = 859.64703
$ws.Range("M107").Value = 1060.35297

$ws.Range("H122").Value = 5295.619
$ws.Range("I122").Value = 5213.625
$ws.Range("K122").Value = 15640.875
$ws.Range("M122").Value = -13190.875

$ws.Range("H136").Value = 1910.4445
$ws.Range("I136").Value = 1242
$ws.Range("K136").Value = 3726
$ws.Range("M136").Value = -1176

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2.909091
$ws.Range("I12").Value = 3.6666667
$ws.Range("K12").Value = 11.0000001
$ws.Range("M12").Value = 161.9999999

$ws.Range("H140").Value = 5960
$ws.Range("I140").Value = 2698
$ws.Range("K140").Value = 8094
$ws.Range("M140").Value = -2914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 30238.072
$ws.Range("J94").Value = 30238.072
$ws.Range("L94").Value = 30238.072
$ws.Range("N94").Value = -31590.072

$ws.Range("H97").Value = 6220
$ws.Range("I97").Value = 643.5714
$ws.Range("K97").Value = 643.5714
$ws.Range("M97").Value = -147.5714

$ws.Range("H107").Value = 759.0833
$ws.Range("I107").Value = 445.4
$ws.Range("J107").Value = 983.1429000000001
$ws.Range("K107").Value = 445.4
$ws.Range("L107").Value = 983.1429000000001
$ws.Range("M107").Value = 1474.6
$ws.Range("N107").Value = -4823.1429

$ws.Range("H122").Value = 2593.3447
$ws.Range("I122").Value = 2465.1052
$ws.Range("J122").Value = 2837
$ws.Range("K122").Value = 7395.3156
$ws.Range("L122").Value = 8511
$ws.Range("M122").Value = -4945.3156
$ws.Range("N122").Value = -13411

$ws.Range("H126").Value = 4351.1113
$ws.Range("I126").Value = 4462.1665
$ws.Range("J126").Value = 4129
$ws.Range("K126").Value = 13386.4995
$ws.Range("L126").Value = 12387
$ws.Range("M126").Value = -10916.4995
$ws.Range("N126").Value = -17327

$ws.Range("H132").Value = 5174.9375
$ws.Range("I132").Value = 5592.7856
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 16778.3568
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -14248.3568
$ws.Range("N132").Value = -11810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5077.95
$ws.Range("I93").Value = 3826.6875
$ws.Range("K93").Value = 3826.6875
$ws.Range("M93").Value = -2578.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 500
$ws.Range("J62").Value = 500
$ws.Range("L62").Value = 500
$ws.Range("N62").Value = -1748

$ws.Range("H65").Value = 500
$ws.Range("J65").Value = 500
$ws.Range("L65").Value = 2500
$ws.Range("N65").Value = -8740

$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630

$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184

$ws.Range("I100").Value = 1525.3334
$ws.Range("J100").Value = 899
$ws.Range("K100").Value = 3050.6668
$ws.Range("L100").Value = 1798
$ws.Range("M100").Value = -2509.6668
$ws.Range("N100").Value = -2880

$ws.Range("H107").Value = 990.13336
$ws.Range("I107").Value = 1120
$ws.Range("J107").Value = 795.3333
$ws.Range("K107").Value = 3360
$ws.Range("L107").Value = 2385.9999
$ws.Range("M107").Value = -1440
$ws.Range("N107").Value = -6225.9999

$ws.Range("H113").Value = 1788.5
$ws.Range("I113").Value = 777
$ws.Range("K113").Value = 2331
$ws.Range("M113").Value = -161
